$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 22910
$ws.Range("E2").Value = 1621
$ws.Range("F2").Value = 1576
$ws.Range("G2").Value = 1335
$ws.Range("H2").Value = 1051
$ws.Range("I2").Value = 1068
$ws.Range("J2").Value = -17
$ws.Range("K2").Value = 22944
$ws.Range("L2").Value = 13009
$ws.Range("M2").Value = 9935
$ws.Range("N2").Value = 9929
$ws.Range("O2").Value = 5
$ws.Range("P2").Value = 1500
$ws.Range("Q2").Value = 1622
$ws.Range("R2").Value = -645
$ws.Range("S2").Value = 429
$ws.Range("T2").Value = 694
$ws.Range("U2").Value = 928
$ws.Range("V2").Value = 7562
$ws.Range("W2").Value = 7.07
$ws.Range("X2").Value = 4.59
$ws.Range("Y2").Value = 11.13
$ws.Range("Z2").Value = 4.78
$ws.Range("AA2").Value = 130.95
$ws.Range("AB2").Value = 578.97
$ws.Range("AC2").Value = 3561
$ws.Range("AD2").Value = 16.8
$ws.Range("AE2").Value = 33860
$ws.Range("AF2").Value = 1.77
$ws.Range("AG2").Value = 1300
$ws.Range("AH2").Value = 2.17
$ws.Range("AI2").Value = 35.69
$ws.Range("AJ2").Value = 30000000

# Row 3
$ws.Range("D3").Value = 22017
$ws.Range("E3").Value = 1544
$ws.Range("F3").Value = 1544
$ws.Range("G3").Value = 1136
$ws.Range("H3").Value = 706
$ws.Range("I3").Value = 703
$ws.Range("J3").Value = 3
$ws.Range("K3").Value = 22517
$ws.Range("L3").Value = 12370
$ws.Range("M3").Value = 10147
$ws.Range("N3").Value = 10133
$ws.Range("O3").Value = 14
$ws.Range("P3").Value = 1500
$ws.Range("Q3").Value = 1003
$ws.Range("R3").Value = -959
$ws.Range("S3").Value = -506
$ws.Range("T3").Value = 690
$ws.Range("U3").Value = 313
$ws.Range("V3").Value = 7402
$ws.Range("W3").Value = 7.01
$ws.Range("X3").Value = 3.21
$ws.Range("Y3").Value = 7.01
$ws.Range("Z3").Value = 3.1
$ws.Range("AA3").Value = 121.91
$ws.Range("AB3").Value = 591.9
$ws.Range("AC3").Value = 2343
$ws.Range("AD3").Value = 19.63
$ws.Range("AE3").Value = 34553
$ws.Range("AF3").Value = 1.33
$ws.Range("AG3").Value = 1000
$ws.Range("AH3").Value = 2.17
$ws.Range("AI3").Value = 41.72
$ws.Range("AJ3").Value = 30000000

# Row 4
$ws.Range("D4").Value = 22136
$ws.Range("E4").Value = 1244
$ws.Range("F4").Value = 1244
$ws.Range("G4").Value = 1020
$ws.Range("H4").Value = 807
$ws.Range("I4").Value = 807
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 22618
$ws.Range("L4").Value = 11957
$ws.Range("M4").Value = 10661
$ws.Range("N4").Value = 10646
$ws.Range("O4").Value = 15
$ws.Range("P4").Value = 1500
$ws.Range("Q4").Value = 2422
$ws.Range("R4").Value = -1020
$ws.Range("S4").Value = -873
$ws.Range("T4").Value = 431
$ws.Range("U4").Value = 1991
$ws.Range("V4").Value = 6824
$ws.Range("W4").Value = 5.62
$ws.Range("X4").Value = 3.65
$ws.Range("Y4").Value = 7.77
$ws.Range("Z4").Value = 3.58
$ws.Range("AA4").Value = 112.16
$ws.Range("AB4").Value = 627.59
$ws.Range("AC4").Value = 2690
$ws.Range("AD4").Value = 14.77
$ws.Range("AE4").Value = 36303
$ws.Range("AF4").Value = 1.09
$ws.Range("AG4").Value = 800
$ws.Range("AH4").Value = 2.01
$ws.Range("AI4").Value = 29.07
$ws.Range("AJ4").Value = 30000000

# Row 5
$ws.Range("D5").Value = 23437
$ws.Range("E5").Value = 1584
$ws.Range("F5").Value = 1584
$ws.Range("G5").Value = 1385
$ws.Range("H5").Value = 1060
$ws.Range("I5").Value = 1051
$ws.Range("J5").Value = 9
$ws.Range("K5").Value = 22558
$ws.Range("L5").Value = 11077
$ws.Range("M5").Value = 11481
$ws.Range("N5").Value = 11457
$ws.Range("O5").Value = 24
$ws.Range("P5").Value = 1500
$ws.Range("Q5").Value = 1446
$ws.Range("R5").Value = -1638
$ws.Range("S5").Value = -541
$ws.Range("T5").Value = 456
$ws.Range("U5").Value = 990
$ws.Range("V5").Value = 6484
$ws.Range("W5").Value = 6.76
$ws.Range("X5").Value = 4.52
$ws.Range("Y5").Value = 9.51
$ws.Range("Z5").Value = 4.69
$ws.Range("AA5").Value = 96.48
$ws.Range("AB5").Value = 685.28
$ws.Range("AC5").Value = 3503
$ws.Range("AD5").Value = 18.58
$ws.Range("AE5").Value = 39070
$ws.Range("AG5").Value = 1100
$ws.Range("AH5").Value = 1.69
$ws.Range("AI5").Value = 30.69
$ws.Range("AJ5").Value = 30000000

# Row 6
$ws.Range("D6").Value = 24850
$ws.Range("E6").Value = 2050
$ws.Range("F6").Value = 2050
$ws.Range("G6").Value = 1778
$ws.Range("H6").Value = 1322
$ws.Range("I6").Value = 1316
$ws.Range("K6").Value = 24058
$ws.Range("L6").Value = 11084
$ws.Range("M6").Value = 12975
$ws.Range("N6").Value = 12945
$ws.Range("P6").Value = 1500
$ws.Range("Q6").Value = 2047
$ws.Range("R6").Value = 673
$ws.Range("S6").Value = -806
$ws.Range("T6").Value = 552
$ws.Range("U6").Value = 1494
$ws.Range("V6").Value = 6102
$ws.Range("W6").Value = 8.25
$ws.Range("X6").Value = 5.32
$ws.Range("Y6").Value = 10.78
$ws.Range("Z6").Value = 5.67
$ws.Range("AA6").Value = 85.42
$ws.Range("AB6").Value = 789.6799999999999
$ws.Range("AC6").Value = 4385
$ws.Range("AD6").Value = 11.17
$ws.Range("AE6").Value = 44144
$ws.Range("AF6").Value = 1.11
$ws.Range("AG6").Value = 1200
$ws.Range("AH6").Value = 2.45
$ws.Range("AI6").Value = 26.75
$ws.Range("AJ6").Value = 30000000

# Row 7
$ws.Range("D7").Value = 23236
$ws.Range("E7").Value = 1767
$ws.Range("G7").Value = 1700
$ws.Range("H7").Value = 1237
$ws.Range("I7").Value = 1243
$ws.Range("K7").Value = 25074
$ws.Range("L7").Value = 11219
$ws.Range("M7").Value = 13855
$ws.Range("N7").Value = 13831
$ws.Range("P7").Value = 1500
$ws.Range("Q7").Value = 1870
$ws.Range("R7").Value = -959
$ws.Range("S7").Value = -448
$ws.Range("T7").Value = 665
$ws.Range("U7").Value = 1953
$ws.Range("W7").Value = 7.61
$ws.Range("X7").Value = 5.33
$ws.Range("Y7").Value = 9.289999999999999
$ws.Range("Z7").Value = 5.04
$ws.Range("AA7").Value = 80.97
$ws.Range("AC7").Value = 4144
$ws.Range("AD7").Value = 13.03
$ws.Range("AE7").Value = 47166
$ws.Range("AF7").Value = 1.14
$ws.Range("AG7").Value = 1225
$ws.Range("AH7").Value = 2.27
$ws.Range("AI7").Value = 29.56

# Row 8
$ws.Range("D8").Value = 24594
$ws.Range("E8").Value = 1927
$ws.Range("G8").Value = 1782
$ws.Range("H8").Value = 1348
$ws.Range("I8").Value = 1342
$ws.Range("K8").Value = 25687
$ws.Range("L8").Value = 10941
$ws.Range("M8").Value = 14746
$ws.Range("N8").Value = 14696
$ws.Range("P8").Value = 1500
$ws.Range("Q8").Value = 1969
$ws.Range("R8").Value = -686
$ws.Range("S8").Value = -601
$ws.Range("T8").Value = 692
$ws.Range("U8").Value = 1348
$ws.Range("W8").Value = 7.83
$ws.Range("X8").Value = 5.48
$ws.Range("Y8").Value = 9.41
$ws.Range("Z8").Value = 5.31
$ws.Range("AA8").Value = 74.2
$ws.Range("AC8").Value = 4472
$ws.Range("AD8").Value = 11.78
$ws.Range("AE8").Value = 50113
$ws.Range("AF8").Value = 1.05
$ws.Range("AG8").Value = 1220
$ws.Range("AH8").Value = 2.31
$ws.Range("AI8").Value = 27.28

# Row 9
$ws.Range("D9").Value = 25718
$ws.Range("E9").Value = 2150
$ws.Range("G9").Value = 2044
$ws.Range("H9").Value = 1538
$ws.Range("I9").Value = 1530
$ws.Range("K9").Value = 26470
$ws.Range("L9").Value = 10585
$ws.Range("M9").Value = 15885
$ws.Range("N9").Value = 15848
$ws.Range("P9").Value = 1500
$ws.Range("Q9").Value = 1836
$ws.Range("R9").Value = -637
$ws.Range("S9").Value = -524
$ws.Range("T9").Value = 685
$ws.Range("U9").Value = 1515
$ws.Range("W9").Value = 8.359999999999999
$ws.Range("X9").Value = 5.98
$ws.Range("Y9").Value = 10.02
$ws.Range("Z9").Value = 5.9
$ws.Range("AA9").Value = 66.63
$ws.Range("AC9").Value = 5102
$ws.Range("AD9").Value = 10.33
$ws.Range("AE9").Value = 54041
$ws.Range("AF9").Value = 0.98
$ws.Range("AG9").Value = 1250
$ws.Range("AH9").Value = 2.37
$ws.Range("AI9").Value = 24.5
